# TestData.xlsx ("Extent Report generated by diff approach"):
# cell AB3 (the "id" column value for test-case row 3 / TC_002) is
# updated from "779" to "481". The value is a text string (not a
# number), so a leading apostrophe is used to force text entry and
# keep it stored as a shared-string cell, matching the original cell's
# text type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AB3").Value = "'481"
